$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the two rows that were dropped entirely from the sheet:
#    "RM 232" (row 26) and "SC 92" (row 28). Deleting row 26 first
#    shifts "SC 92" up to row 27, so delete row 27 next.
# ------------------------------------------------------------------
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# ------------------------------------------------------------------
# 2) Apply the remaining cell-level edits (using the row numbers as
#    they stand after the two deletions above).
# ------------------------------------------------------------------

# RM 14 (row 5): F5 value removed
$ws.Range("F5").ClearContents()

# RM 21 (row 6): D6 gains a value
$ws.Range("D6").Value = -14.2

# RM 38 (row 8): D8 value removed
$ws.Range("D8").ClearContents()

# RM 58 (row 11): F11 gains a value
$ws.Range("F11").Value = 17.65

# RM 125 (row 19): D19 gains a value, F19 value removed
$ws.Range("D19").Value = -15.5
$ws.Range("F19").ClearContents()

# RM 135 (row 21): D21 value removed
$ws.Range("D21").ClearContents()

# RM 140 (row 23): D23 and F23 gain values
$ws.Range("D23").Value = -13.9
$ws.Range("F23").Value = 16.48

# RM 145 (row 25): F25 gains a value
$ws.Range("F25").Value = 16.6

# SC 5 (row 26): B26 value removed
$ws.Range("B26").ClearContents()

# SC 101 (row 27): B27 gains a value, D27 and F27 values removed
$ws.Range("B27").Value = -20.4
$ws.Range("D27").ClearContents()
$ws.Range("F27").ClearContents()

# SC 119 (row 29): B29 value removed, D29 gains a value, F29 value removed
$ws.Range("B29").ClearContents()
$ws.Range("D29").Value = -13.0
$ws.Range("F29").ClearContents()

# SC 120 (row 30): F30 gains a value
$ws.Range("F30").Value = 16.89

# SC 232 (row 33): F33 gains a value
$ws.Range("F33").Value = 17.53
